# Update countries & provincias Spain
# Applies the data refresh described by the commit: updated case numbers for
# several countries, a new country (Georgia) overtaking three others in the
# "Casos totales" ranking, a ranking swap between Santa Lucia / Timor Oriental,
# and a refreshed "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Septiembre de 2020 a las 09:26"

# 2) Row 64 - Armenia: updated case counts
$ws.Range("B64").Value = 49072
$ws.Range("C64").Value = 429
$ws.Range("D64").Value = 43576
$ws.Range("E64").Value = 4548
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 948

# 3) Row 69 - Afganistan: updated case counts
$ws.Range("B69").Value = 39192
$ws.Range("C69").Value = 6
$ws.Range("D69").Value = 32635
$ws.Range("E69").Value = 5104
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 1453

# 4) Row 78 - Australia: updated case counts
$ws.Range("B78").Value = 27015
$ws.Range("C78").Value = 15
$ws.Range("D78").Value = 24571
$ws.Range("E78").Value = 1574
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 870

# 5) Row 82 - Hungria: updated case counts
$ws.Range("B82").Value = 23077
$ws.Range("C82").Value = 950
$ws.Range("D82").Value = 5099
$ws.Range("E82").Value = 17248
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 12
$ws.Range("H82").Value = 730

# 6) Rows 126-129: Georgia's case total (4960) now overtakes Surinam (4817),
#    Republica de Africa Central (4806) and Ruanda (4798), so Georgia's row
#    moves up to 126 and the other three each shift down one row.
$ws.Range("A126").Value = "Georgia"
$ws.Range("B126").Value = 4960
$ws.Range("C126").Value = 296
$ws.Range("D126").Value = 1819
$ws.Range("E126").Value = 3114
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 27

$ws.Range("A127").Value = "Surinam"
$ws.Range("B127").Value = 4817
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 4596
$ws.Range("E127").Value = 119
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 102

$ws.Range("A128").Value = "Republica de Africa Central"
$ws.Range("B128").Value = 4806
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 1840
$ws.Range("E128").Value = 2904
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 62

$ws.Range("A129").Value = "Ruanda"
$ws.Range("B129").Value = 4798
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 3080
$ws.Range("E129").Value = 1689
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 29

# 7) Rows 205-206: Timor Oriental and Santa Lucia are tied on total cases (27)
#    and swap rank order; the underlying figures are identical for both.
$ws.Range("A205").Value = "Timor Oriental"
$ws.Range("B205").Value = 27
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 27
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

$ws.Range("A206").Value = "Santa Lucia"
$ws.Range("B206").Value = 27
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 27
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
